$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is a plain number-like string must be forced to Text
# format first, otherwise Excel auto-converts them to numeric values and
# mangles formatting (trailing zeros, thousands separators used as decimal
# grouping in these European-style prices, etc.)

# Row 2
$ws.Range("D2").Value = "62.117.58"
$ws.Range("E2").Value = "  +2.72%  "

# Row 3
$ws.Range("D3").Value = "2.438.07"
$ws.Range("E3").Value = "  +4.57%  "

# Row 4
$ws.Range("E4").Value = "  +0.04%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "557.00"
$ws.Range("E5").Value = "  +2.10%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "138.65"

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  +0.02%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.584"
$ws.Range("E8").Value = "  +1.10%  "

# Row 9
$ws.Range("D9").Value = "2.435.14"
$ws.Range("E9").Value = "  +4.54%  "

# Row 10
$ws.Range("E10").Value = "  +2.54%  "

# Row 11
$ws.Range("E11").Value = "  +4.20%  "

# Row 12
$ws.Range("E12").Value = "  +0.04%  "

# Row 13
$ws.Range("E13").Value = "  +4.04%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "25.86"
$ws.Range("E14").Value = "  +9.40%  "

# Row 15
$ws.Range("D15").Value = "2.873.54"
$ws.Range("E15").Value = "  +4.61%  "

# Row 16
$ws.Range("D16").Value = "62.049.11"
$ws.Range("E16").Value = "  +2.72%  "

# Row 17
$ws.Range("E17").Value = "  +5.64%  "

# Row 18
$ws.Range("D18").Value = "2.444.96"
$ws.Range("E18").Value = "  +5.25%  "

# Row 19
$ws.Range("E19").Value = "  +5.16%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "345.80"
$ws.Range("E20").Value = "  +9.89%  "

# Row 21
$ws.Range("E21").Value = "  +2.42%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.84"
$ws.Range("E22").Value = "  +2.89%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "65.16"
$ws.Range("E24").Value = "  +1.67%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.173"
$ws.Range("E25").Value = "  +0.71%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.999"
$ws.Range("E26").Value = "  -0.14%  "

# Row 27
$ws.Range("E27").Value = "  +11.10%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.28"
$ws.Range("E28").Value = "  +5.75%  "

# Row 29
$ws.Range("E29").Value = "  +13.04%  "

# Row 30
$ws.Range("D30").Value = "0.0₃0788"
$ws.Range("E30").Value = "  +7.53%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.81"
$ws.Range("E31").Value = "  +4.80%  "

# Row 32
$ws.Range("E32").Value = "  +6.67%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "171.18"
$ws.Range("E33").Value = "  -0.86%  "

# Row 34
$ws.Range("E34").Value = "  +5.60%  "

# Row 35
$ws.Range("E35").Value = "  +4.27%  "

# Row 36
$ws.Range("B36").Value = "EthereumClassic"
$ws.Range("C36").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "18.58"
$ws.Range("E36").Value = "  +4.14%  "

# Row 37
$ws.Range("B37").Value = "Bittensor"
$ws.Range("C37").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "374.50"
$ws.Range("E37").Value = "  +16.39%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.48"

# Row 39
$ws.Range("E39").Value = "  -0.02%  "

# Row 40
$ws.Range("E40").Value = "  -0.05%  "

# Row 41
$ws.Range("E41").Value = "  +10.25%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "39.28"
$ws.Range("E42").Value = "  +3.60%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "146.61"
$ws.Range("E43").Value = "  +6.72%  "

# Row 44
$ws.Range("E44").Value = "  +5.34%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "20.73"
$ws.Range("E45").Value = "  +8.53%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0959"
$ws.Range("E46").Value = "  +2.04%  "

# Row 47
$ws.Range("E47").Value = "  +4.47%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0519"
$ws.Range("E48").Value = "  +5.04%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "18.02"
$ws.Range("E49").Value = "  +6.64%  "

# Row 50
$ws.Range("E50").Value = "  +4.03%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.71"
$ws.Range("E51").Value = "  +11.49%  "
